$d = $word.ActiveDocument
$word.Options.ShowHiddenBookmarks = $true

# Locate the target list paragraph.
$para = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $cand = $d.Paragraphs.Item($i)
    if ($cand.Range.Text -like "*Calculos*duty*cycle*") {
        $para = $cand
        break
    }
}
$pStart = $para.Range.Start

# Goal: turn
#   "...del PWM y de los valores para duty cycle"
# into
#   "...del PWM y de los valores/precisión para duty cycle"
# i.e. insert "/precisión" right after "valores" (before " para ").

# Insert the new text next to "duty" first -- that run has no leading/
# trailing whitespace (so no xml:space="preserve"), which keeps the new
# "/precisión" run clean instead of inheriting xml:space="preserve" from
# the " del PWM y de los valores para " run it would otherwise land in.
$full = $para.Range.Text
$dutyEndPos = $pStart + $full.IndexOf("duty") + 4
$rInsert = $d.Range($dutyEndPos, $dutyEndPos)
$rInsert.InsertBefore("/precisión")

# Split "duty/precisión" into separate "duty" and "/precisión" runs by
# dropping a transient bookmark at the boundary (add+delete forces a run
# split there without Word silently re-merging the two runs).
$d.Bookmarks.Add("tempSplit1", $d.Range($dutyEndPos, $dutyEndPos))
$d.Bookmarks.Item("tempSplit1").Delete()

# Cut the now-standalone "/precisión" run and paste it where it belongs,
# right after "valores".
$rCut = $d.Range($dutyEndPos, $dutyEndPos + 10)
$rCut.Cut()

$full2 = $para.Range.Text
$targetPos = $pStart + $full2.IndexOf("valores") + 7
$rTarget = $d.Range($targetPos, $targetPos)
$rTarget.Paste()

# Word also leaves its "_GoBack" (last-edit) bookmark sitting inside
# "cycle" (right after the "c") when this edit is saved, which splits
# that run in two. Reproduce it.
$full3 = $para.Range.Text
$cyclePos = $pStart + $full3.IndexOf("cycle") + 1
$rGoBack = $d.Range($cyclePos, $cyclePos)
$d.Bookmarks.Add("_GoBack", $rGoBack)

Write-Output $para.Range.Text
